# Auto-generated edit script: updates cryptos list price/volume data
# per commit "Updated cryptos list on Wed Jan 24 20:59:09 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'39.671.84"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = "'2.201.87"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'291.64"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').Value = "'86.29"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.28%  '
$ws.Range('E7').Value = '  +0.92%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('D10').Value = "'30.13"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.91%  '
$ws.Range('E11').Value = '  +1.98%  '
$ws.Range('D12').Value = "'47.39"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').Value = "'6.31"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('D15').Value = "'2.544.20"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Value = "'2.213.04"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('E18').Value = '  +2.47%  '
$ws.Range('D19').Value = "'39.591.50"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').Value = "'11.27"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.13%  '
$ws.Range('D22').Value = "'5.76"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').Value = "'65.50"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('D24').Value = "'235.09"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.96%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  +2.54%  '
$ws.Range('E27').Value = '  +1.79%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').Value = "'9.24"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.25%  '
$ws.Range('D31').Value = "'151.80"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('D32').Value = "'32.45"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.81%  '
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('E35').Value = '  +3.67%  '
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('D37').Value = "'0.112"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('E38').Value = '  +6.16%  '
$ws.Range('D39').Value = "'15.90"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.44%  '
$ws.Range('E40').Value = '  +3.51%  '
$ws.Range('E41').Value = '  +3.26%  '
$ws.Range('D42').Value = "'2.061.89"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.75%  '
$ws.Range('E43').Value = '  +4.83%  '
$ws.Range('D44').Value = "'2.13"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.44%  '
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'9.90"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.44%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'17.58"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.55%  '
$ws.Range('D48').Value = "'2.59"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').Value = "'2.422.56"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').Value = "'70.62"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('D51').Value = "'88.55"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.88%  '
